$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "Total" label and a SUM formula, leaving row 9 empty (matches target layout)
$ws.Range("A10").Value = "Total"
$ws.Range("B10").Formula = "=SUM(B2:B8)"
$ws.Range("B10").Style = "Currency"

# Resize column B to fit the new content
$ws.Columns.Item(2).EntireColumn.AutoFit()

# Update the active selection to A8 (matches the recorded cursor position in the diff)
$ws.Range("A8").Select()
